# Update the 'want-to-go count' (F column) figures across all sheets
# to match the latest scrape, per commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 791
$ws.Range("F3").Value = 14510
$ws.Range("F4").Value = 14722
$ws.Range("F7").Value = 1418
$ws.Range("F8").Value = 5977
$ws.Range("F10").Value = 583
$ws.Range("F11").Value = 114
$ws.Range("F12").Value = 376
$ws.Range("F13").Value = 203
$ws.Range("F14").Value = 1578
$ws.Range("F15").Value = 469
$ws.Range("F16").Value = 2138
$ws.Range("F17").Value = 1251
$ws.Range("F18").Value = 1897
$ws.Range("F21").Value = 2316
$ws.Range("F22").Value = 587
$ws.Range("F23").Value = 845
$ws.Range("F24").Value = 3443
$ws.Range("F25").Value = 311
$ws.Range("F27").Value = 2516
$ws.Range("F28").Value = 627
$ws.Range("F30").Value = 1343
$ws.Range("F31").Value = 1852
$ws.Range("F33").Value = 1494
$ws.Range("F34").Value = 156
$ws.Range("F36").Value = 5089
$ws.Range("F37").Value = 5003
$ws.Range("F38").Value = 316
$ws.Range("F39").Value = 698
$ws.Range("F44").Value = 352
$ws.Range("F45").Value = 132
$ws.Range("F46").Value = 104
$ws.Range("F47").Value = 4465
$ws.Range("F48").Value = 648
$ws.Range("F49").Value = 313

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 131
$ws.Range("F20").Value = 21

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7792
$ws.Range("F3").Value = 280
$ws.Range("F4").Value = 968

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7792
$ws.Range("F3").Value = 791
$ws.Range("F4").Value = 280
$ws.Range("F5").Value = 968
$ws.Range("F7").Value = 14510
$ws.Range("F8").Value = 14722
$ws.Range("F10").Value = 1418
$ws.Range("F11").Value = 5977
$ws.Range("F13").Value = 131
$ws.Range("F16").Value = 1578
$ws.Range("F17").Value = 469
$ws.Range("F19").Value = 845
$ws.Range("F21").Value = 3443
$ws.Range("F23").Value = 2516
$ws.Range("F24").Value = 627
$ws.Range("F26").Value = 1852
$ws.Range("F34").Value = 21
$ws.Range("F35").Value = 5089
$ws.Range("F36").Value = 5003
$ws.Range("F37").Value = 316
$ws.Range("F38").Value = 698
$ws.Range("F41").Value = 352
$ws.Range("F42").Value = 132
$ws.Range("F44").Value = 104
$ws.Range("F45").Value = 648
$ws.Range("F46").Value = 313
